$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.173.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.09%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.576.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.65%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.57%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'208.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.36%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -3.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'0.0610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.47%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.94%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.00%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.48%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.798.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.64%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.592.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'Polkadot"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'4.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.13%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.04%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.93%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.166.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.07%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -1.85%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.55%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'209.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.51%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.93%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -2.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.21%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'143.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.05%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.54%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.44%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.89%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.42%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.76%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.50%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.280.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.41%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.55%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.607"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.51%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.72%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.10%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'WEMIXToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -10.72%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.64%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.46%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +2.92%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'MXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'2.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.74%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'TrustWalletToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.764"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.84%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'62.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.711.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.61%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'88.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.72%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -2.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -4.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.44%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.51%  "
$ws.Range("E51").Style = "Normal"
